# Update "想去人数" (column F) values on the "展览" and "全部类型" sheets
# to reflect newly scraped counts, per commit "Update gh-pages to output
# generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# Row -> new F value for worksheet "展览"
$sheetExhibition = @{
    2  = 638
    3  = 2229
    5  = 13379
    6  = 80
    7  = 122
    8  = 524
    9  = 490
    10 = 1198
    11 = 1009
    12 = 13825
    13 = 14510
    17 = 44
    21 = 9
    22 = 45
    24 = 1110
    26 = 59
    27 = 5556
    29 = 1038
    30 = 5354
    31 = 31
    32 = 20
    33 = 137
}

# Row -> new F value for worksheet "全部类型"
$sheetAllTypes = @{
    2  = 638
    3  = 2229
    5  = 13379
    6  = 80
    8  = 122
    9  = 524
    10 = 490
    11 = 1198
    12 = 1009
    13 = 13825
    14 = 14510
    18 = 44
    22 = 9
    23 = 45
    25 = 1110
    27 = 59
    28 = 5556
    30 = 1038
    31 = 5354
    32 = 31
    33 = 20
    34 = 137
}

$ws1 = $wb.Worksheets.Item("展览")
foreach ($row in $sheetExhibition.Keys) {
    $ws1.Cells.Item($row, 6).Value = $sheetExhibition[$row]
}

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($row in $sheetAllTypes.Keys) {
    $ws4.Cells.Item($row, 6).Value = $sheetAllTypes[$row]
}
